# Update the "想去人数" (F column) counts across the workbook's sheets
# to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 558
$ws1.Range("F3").Value  = 10388
$ws1.Range("F6").Value  = 6941
$ws1.Range("F8").Value  = 134
$ws1.Range("F9").Value  = 12091
$ws1.Range("F10").Value = 12690
$ws1.Range("F11").Value = 1309
$ws1.Range("F12").Value = 1279
$ws1.Range("F13").Value = 5397
$ws1.Range("F19").Value = 337
$ws1.Range("F20").Value = 1999
$ws1.Range("F23").Value = 879
$ws1.Range("F29").Value = 2015
$ws1.Range("F32").Value = 993
$ws1.Range("F34").Value = 37
$ws1.Range("F37").Value = 4370
$ws1.Range("F40").Value = 591
$ws1.Range("F41").Value = 642
$ws1.Range("F42").Value = 687
$ws1.Range("F46").Value = 25

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 69
$ws2.Range("F15").Value = 5
$ws2.Range("F18").Value = 6

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6451

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 558
$ws4.Range("F3").Value  = 228
$ws4.Range("F9").Value  = 12091
$ws4.Range("F10").Value = 12690
$ws4.Range("F12").Value = 1309
$ws4.Range("F13").Value = 1279
$ws4.Range("F28").Value = 2015
$ws4.Range("F31").Value = 11
$ws4.Range("F37").Value = 4370
$ws4.Range("F40").Value = 125
$ws4.Range("F41").Value = 642
$ws4.Range("F42").Value = 687
$ws4.Range("F48").Value = 0
